$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D2").Value = "25.684.05"
$ws.Range("E2").Value = "  -3.62%  "
$ws.Range("D3").Value = "1.742.88"
$ws.Range("E3").Value = "  -5.53%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "238.00"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "0.4944"
$ws.Range("E7").Value = "  -6.34%  "
$ws.Range("D8").Value = "41.58"
$ws.Range("E8").Value = "  -7.60%  "
$ws.Range("D9").Value = "0.2388"
$ws.Range("E9").Value = "  -24.34%  "
$ws.Range("E10").Value = "  -12.45%  "
$ws.Range("D11").Value = "1.740.47"
$ws.Range("E11").Value = "  -5.75%  "
$ws.Range("D12").Value = "0.06843"
$ws.Range("E12").Value = "  -12.05%  "
$ws.Range("D13").Value = "14.61"
$ws.Range("E13").Value = "  -23.26%  "
$ws.Range("E14").Value = "  -11.29%  "
$ws.Range("D15").Value = "77.18"
$ws.Range("E15").Value = "  -12.65%  "
$ws.Range("D16").Value = "0.5785"
$ws.Range("E16").Value = "  -26.43%  "
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").Value = "25.736.14"
$ws.Range("E19").Value = "  -3.50%  "
$ws.Range("D20").Value = "11.45"
$ws.Range("E20").Value = "  -17.72%  "
$ws.Range("D21").Value = "0.000006445"
$ws.Range("E21").Value = "  -18.72%  "
$ws.Range("D22").Value = "1.960.58"
$ws.Range("E22").Value = "  -6.12%  "
$ws.Range("D23").Value = "3.959"
$ws.Range("E23").Value = "  -14.11%  "
$ws.Range("D24").Value = "5.052"
$ws.Range("E24").Value = "  -15.74%  "
$ws.Range("D25").Value = "7.785"
$ws.Range("E25").Value = "  -16.81%  "
$ws.Range("D26").Value = "136.60"
$ws.Range("E26").Value = "  -4.40%  "
$ws.Range("D27").Value = "1.472"
$ws.Range("E27").Value = "  -12.51%  "
$ws.Range("D28").Value = "1.829"
$ws.Range("E28").Value = "  -17.89%  "
$ws.Range("D29").Value = "14.52"
$ws.Range("E29").Value = "  -14.80%  "
$ws.Range("D30").Value = "100.71"
$ws.Range("E30").Value = "  -9.30%  "
$ws.Range("D31").Value = "3.796"
$ws.Range("E31").Value = "  -9.79%  "
$ws.Range("D32").Value = "0.08122"
$ws.Range("E32").Value = "  -6.70%  "
$ws.Range("D33").Value = "3.344"
$ws.Range("E33").Value = "  -18.24%  "
$ws.Range("D34").Value = "0.04359"
$ws.Range("E34").Value = "  -10.77%  "
$ws.Range("D35").Value = "1.001"
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("D36").Value = "2.685"
$ws.Range("E36").Value = "  -6.19%  "
$ws.Range("D37").Value = "1.021"
$ws.Range("E37").Value = "  -10.65%  "
$ws.Range("D38").Value = "0.6073"
$ws.Range("D39").Value = "2.695"
$ws.Range("E39").Value = "  -13.28%  "
$ws.Range("D40").Value = "2.079"
$ws.Range("E40").Value = "  -9.79%  "
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").Value = "103.12"
$ws.Range("E42").Value = "  -6.16%  "
$ws.Range("D43").Value = "0.01483"
$ws.Range("E43").Value = "  -14.44%  "
$ws.Range("D44").Value = "0.7797"
$ws.Range("E44").Value = "  -13.52%  "
$ws.Range("D45").Value = "5.133"
$ws.Range("E45").Value = "  -13.59%  "
$ws.Range("D46").Value = "0.3765"
$ws.Range("E46").Value = "  -22.05%  "
$ws.Range("E47").Value = "  -12.39%  "
$ws.Range("D48").Value = "5.980"
$ws.Range("E48").Value = "  -22.51%  "
$ws.Range("D49").Value = "0.1062"
$ws.Range("E49").Value = "  -14.73%  "
$ws.Range("D50").Value = "30.09"
$ws.Range("E50").Value = "  -13.66%  "
$ws.Range("D51").Value = "52.62"
$ws.Range("E51").Value = "  -12.29%  "
